# Remove the post entry that was deleted from the spreadsheet
# (「幸せは…全然思いもよらない所からやって来る」 …), which lived in row 633.
# Deleting the entire row shifts every subsequent row (634-677) up by one,
# matching the new dimension A1:C676.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(633).Delete()
